$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# Output section originally contains 3 example paragraphs (each framed with a
# bottom paragraph border):
#   9:  "Et lite antall ..." (also carries a stray _GoBack bookmark)
#   11: "Et stort antall ... variert rekke formater ..."
#   13: "Et stort antall ... DOC format ..."
# The edit:
#   - drops the bottom-border formatting from every example paragraph
#   - drops the stray _GoBack bookmark
#   - turns paragraph 11 into a lone "AND/OR" separator paragraph, and
#     re-inserts its former text (+ another "AND/OR" separator) as brand new
#     paragraphs right after it, each wrapped in its own blank paragraph.
# ---------------------------------------------------------------------------

$stortVariertText = "Et stort antall konverteringer er registrert, hvorav alle er gjort til PDF format fra en variert rekke formater. Ingenting peker på en nødvendighet for videre testing."

# 1) Insert the new paragraphs before paragraph 13 (the DOC-format example):
#      <stort/variert text>
#      <blank>
#      AND/OR
#      <blank>
#    This duplicates the text removed from paragraph 11 below, plus adds a
#    second "AND/OR" separator before the final example.
$insertionPoint = $d.Paragraphs.Item(13).Range.Duplicate
$insertionPoint.Collapse(1)
$insertionPoint.InsertBefore("$stortVariertText`r`rAND/OR`r`r")

# 2) Turn (old) paragraph 11 into the first "AND/OR" separator paragraph.
$p11 = $d.Paragraphs.Item(11)
$r11 = $p11.Range
$oldText11 = $r11.Text.TrimEnd([char]13)
$r11.Find.Execute($oldText11, $true, $false, $false, $false, $false, $true, 1, $false, "AND/OR", 2)

# 3) Strip the bottom paragraph border from every example paragraph. The
#    four brand-new paragraphs (13-16) all inherited the border from their
#    insertion point (old paragraph 13), including the two blank spacers
#    among them, so clear all of them too.
$d.Paragraphs.Item(9).Borders.Item(-3).LineStyle = 0
$d.Paragraphs.Item(11).Borders.Item(-3).LineStyle = 0
$d.Paragraphs.Item(13).Borders.Item(-3).LineStyle = 0
$d.Paragraphs.Item(14).Borders.Item(-3).LineStyle = 0
$d.Paragraphs.Item(15).Borders.Item(-3).LineStyle = 0
$d.Paragraphs.Item(16).Borders.Item(-3).LineStyle = 0
$d.Paragraphs.Item(17).Borders.Item(-3).LineStyle = 0

# 4) Remove the stray _GoBack bookmark left over in paragraph 9.
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks.Item("_GoBack").Delete()
}
